$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab (and workbook.xml sheet name) from the 04-18 snapshot to 04-19
$ws.Name = "Through 2022-04-19"

# Update the "2022 (through 04-18)" column header shared string to 04-19
$ws.Range("I1").Value = "2022 (through 04-19)"

# Update April's 2022-to-date carjacking count (row 5 = April)
$ws.Range("I5").Value = 84

# Update the Total row's 2022-to-date carjacking count (row 14 = Total)
$ws.Range("I14").Value = 519
